$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 93500
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = 94121
$ws.Range("B2").Style = "Normal"
$ws.Range("E2").Value = 53
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = 'Vedtrappmossa'
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = 'Crossocalyx hellerianus'
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = '(Nees ex Lindenb.) Meyl.'
$ws.Range("H2").Style = "Normal"
$ws.Range("P2").Value = 'Bäcken vid Mellantjärnsbodarna, Mpd'
$ws.Range("P2").Style = "Normal"
$ws.Range("Q2").Value = 540775.1640602688
$ws.Range("Q2").Style = "Normal"
$ws.Range("R2").Value = 6941931.758068252
$ws.Range("R2").Style = "Normal"
$ws.Range("S2").Value = 10
$ws.Range("S2").Style = "Normal"
$ws.Range("Y2").Value = '''2009-09-02'
$ws.Range("Y2").Style = "Normal"
$ws.Range("AA2").Value = '''2009-09-02'
$ws.Range("AA2").Style = "Normal"
$ws.Range("AC2").Value = 'Jonas Salmonsson'
$ws.Range("AC2").Style = "Normal"
$ws.Range("AH2").ClearContents()
$ws.Range("AI2").ClearContents()
$ws.Range("AN2").Value = 2
$ws.Range("AN2").Style = "Normal"
$ws.Range("AO2").Value = '2 substratenheter # Timmer'
$ws.Range("AO2").Style = "Normal"
$ws.Range("AR2").ClearContents()
$ws.Range("AX2").Value = 'Via Andreas Karlberg'
$ws.Range("AX2").Style = "Normal"

# --- Row 4 ---
$ws.Range("A4").Value = 1866240
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = 73678
$ws.Range("B4").Style = "Normal"
$ws.Range("D4").Value = 'LC'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = 6439
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = 'Gulnål'
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = 'Chaenotheca brachypoda'
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = '(Ach.) Tibell'
$ws.Range("H4").Style = "Normal"
$ws.Range("Q4").Value = 540827.0929000516
$ws.Range("Q4").Style = "Normal"
$ws.Range("R4").Value = 6941759.694429157
$ws.Range("R4").Style = "Normal"
$ws.Range("AH4").Value = 'Granskog'
$ws.Range("AH4").Style = "Normal"

# --- Row 5 ---
$ws.Range("A5").Value = 168997
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = 73685
$ws.Range("B5").Style = "Normal"
$ws.Range("D5").Value = 'VU'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = 492
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = 'Smalskaftslav'
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = 'Chaenotheca gracilenta'
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = '(Ach.) J.Mattsson & Middelb.'
$ws.Range("H5").Style = "Normal"
$ws.Range("P5").Value = 'Täljeån, Mpd'
$ws.Range("P5").Style = "Normal"
$ws.Range("Q5").Value = 540827.0929000516
$ws.Range("Q5").Style = "Normal"
$ws.Range("R5").Value = 6941759.694429157
$ws.Range("R5").Style = "Normal"
$ws.Range("S5").Value = 25
$ws.Range("S5").Style = "Normal"
$ws.Range("Y5").Value = '''2007-11-09'
$ws.Range("Y5").Style = "Normal"
$ws.Range("AA5").Value = '''2007-11-09'
$ws.Range("AA5").Style = "Normal"
$ws.Range("AC5").ClearContents()
$ws.Range("AH5").Value = 'Granskog'
$ws.Range("AH5").Style = "Normal"
$ws.Range("AN5").ClearContents()
$ws.Range("AO5").ClearContents()
$ws.Range("AR5").ClearContents()
$ws.Range("AX5").Value = 'Hans Sundström'
$ws.Range("AX5").Style = "Normal"

# --- Row 6 ---
$ws.Range("A6").Value = 292730
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = 79433
$ws.Range("B6").Style = "Normal"
$ws.Range("D6").Value = 'NT'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = 1049
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = 'Kortskaftad ärgspik'
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = 'Microcalicium ahlneri'
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = 'Tibell'
$ws.Range("H6").Style = "Normal"
$ws.Range("Q6").Value = 540844.6604352774
$ws.Range("Q6").Style = "Normal"
$ws.Range("R6").Value = 6941718.921423005
$ws.Range("R6").Style = "Normal"

# --- Row 7 ---
$ws.Range("A7").Value = 1901472
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = 78569
$ws.Range("B7").Style = "Normal"
$ws.Range("D7").Value = 'NT'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = 6458
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = 'Lunglav'
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = 'Lobaria pulmonaria'
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = '(L.) Hoffm.'
$ws.Range("H7").Style = "Normal"
$ws.Range("Q7").Value = 540812.4636330464
$ws.Range("Q7").Style = "Normal"
$ws.Range("R7").Value = 6941674.310578943
$ws.Range("R7").Style = "Normal"
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()
$ws.Range("AO7").ClearContents()

# --- Row 8 ---
$ws.Range("A8").Value = 1866241
$ws.Range("A8").Style = "Normal"
$ws.Range("Q8").Value = 540835.5837245358
$ws.Range("Q8").Style = "Normal"
$ws.Range("R8").Value = 6941669.529265426
$ws.Range("R8").Style = "Normal"
$ws.Range("AI8").Value = 'Gransumpskog'
$ws.Range("AI8").Style = "Normal"

# --- Row 9 ---
$ws.Range("A9").Value = 1672419
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = 89356
$ws.Range("B9").Style = "Normal"
$ws.Range("D9").Value = 'LC'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = 5447
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = 'Vedticka'
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = 'Fuscoporia viticola'
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value = '(Schwein.) Murrill'
$ws.Range("H9").Style = "Normal"
$ws.Range("Q9").Value = 540838.9380165208
$ws.Range("Q9").Style = "Normal"
$ws.Range("R9").Value = 6941696.743931145
$ws.Range("R9").Style = "Normal"
$ws.Range("AI9").Value = 'Gransumpskog'
$ws.Range("AI9").Style = "Normal"
$ws.Range("AO9").Value = 'Granlåga'
$ws.Range("AO9").Style = "Normal"
